$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime (D) and Correspond Handback DateTime (G)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-22 14:32:35"
$wsZh.Range("D3").Value = "2016-02-22 14:32:35"
$wsZh.Range("G2").Value = "2016-02-22 14:33:25"
$wsZh.Range("G3").Value = "2016-02-22 14:33:25"

# de-de sheet: update Correspond Handoff Datetime (D) and Correspond Handback DateTime (G)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-22 14:32:48"
$wsDe.Range("D3").Value = "2016-02-22 14:32:48"
$wsDe.Range("G2").Value = "2016-02-22 14:33:48"
$wsDe.Range("G3").Value = "2016-02-22 14:33:48"
